$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.545
$ws.Range("B21").Value = 9.658000000000001
$ws.Range("B23").Value = 7.891999999999999
$ws.Range("B25").Value = 6.396
$ws.Range("C27").Value = -13.055
$ws.Range("C31").Value = -13.002
$ws.Range("C39").Value = -12.758
$ws.Range("C48").Value = -11.1
$ws.Range("C51").Value = -11.152
$ws.Range("C52").Value = -11.273
$ws.Range("B53").Value = 6.139999999999999
$ws.Range("C55").Value = -13.368
$ws.Range("C56").Value = -13.222
$ws.Range("B57").Value = 5.575
$ws.Range("C57").Value = -13.385
$ws.Range("B59").Value = 4.678
$ws.Range("B69").Value = 5.667
$ws.Range("C73").Value = -12.575
$ws.Range("B79").Value = 5.771000000000001
$ws.Range("B83").Value = 5.702
$ws.Range("C89").Value = -10.863
$ws.Range("C90").Value = -12.91
$ws.Range("B93").Value = 5.659000000000001
